$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Julio de 2020 a las 21:05"

# --- Swap country names whose table rows traded places ---
# Sudan del Sur / Congo swapped rank (row 115 / 116)
$ws.Range("A115").Value = "Sudan del Sur"
$ws.Range("A116").Value = "Congo"

# Islas Malvinas / Groenlandia swapped rank (row 209 / 210)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Updated statistics per country row ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 3446756
$ws.Range("C4").Value = 32761
$ws.Range("D4").Value = 1535368
$ws.Range("E4").Value = 1773402
$ws.Range("G4").Value = 204
$ws.Range("H4").Value = 137986

# Row 10: España
$ws.Range("B10").Value = 303033
$ws.Range("C10").Value = 681
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 28406

# Row 19: Alemania
$ws.Range("B19").Value = 200378
$ws.Range("C19").Value = 428
$ws.Range("E19").Value = 6141

# Row 21: Francia
$ws.Range("B21").Value = 172377
$ws.Range("C21").Value = 288
$ws.Range("E21").Value = 63960
$ws.Range("G21").Value = 18
$ws.Range("H21").Value = 30029

# Row 30: Suecia
$ws.Range("B30").Value = 75826
$ws.Range("C30").Value = 31
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = 5536

# Row 38: Emiratos Arabes Unidos
$ws.Range("B38").Value = 55198
$ws.Range("C38").Value = 344
$ws.Range("D38").Value = 45513
$ws.Range("E38").Value = 9351
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 334

# Row 68: Uzbekistan
$ws.Range("E68").Value = 5444
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 64

# Row 97: Luxemburgo
$ws.Range("D97").Value = 4183
$ws.Range("E97").Value = 662

# Row 110: Sri Lanka
$ws.Range("B110").Value = 2646
$ws.Range("C110").Value = 29
$ws.Range("E110").Value = 654

# Row 115 (now Sudan del Sur) stats
$ws.Range("B115").Value = 2148
$ws.Range("C115").Value = 127
$ws.Range("D115").Value = 333
$ws.Range("E115").Value = 1774
$ws.Range("G115").Value = 3
$ws.Range("H115").Value = 41

# Row 116 (now Congo) stats
$ws.Range("B116").Value = 2028
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 589
$ws.Range("E116").Value = 1392
$ws.Range("H116").Value = 47
